# Cambio y correcion en numero de cuentas por cobrar cliente
#
# The "Cuenta Contable" numbers used for the client's receivable account
# (111, 111.1, 111.2) were corrected to (113, 113.1, 113.2).
# These three cells are plain inputs; every other cell that shows 111/111.1/111.2
# in the workbook is a formula that references F3/F4/F5 (directly or indirectly),
# so updating the three source cells lets Excel recalculate the rest automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 113
$ws.Range("F4").Value = 113.1
$ws.Range("F5").Value = 113.2

# Recalculate the workbook so every dependent formula (C9, A22, A23, A27,
# A30, A33, K3, K4, ...) picks up the new account numbers.
$excel.CalculateFullRebuild()

# Leave the selection on F4, matching where the edit was made.
$ws.Range("F4").Select()
